# paises.xlsx -- refresh country COVID-19 stats + provincias Spain data pull
# (commit: "Update countries & provincias Spain")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1, col A)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 9 de Julio de 2020 a las 01:35"

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 3156437   # Casos totales
$ws.Cells.Item(4,3).Value = 59353   # Nuevos casos
$ws.Cells.Item(4,4).Value = 1385929   # Casos activos
$ws.Cells.Item(4,5).Value = 1635679   # Recuperados
$ws.Cells.Item(4,7).Value = 857   # Muertes hoy
$ws.Cells.Item(4,8).Value = 134829   # Muertes

# Row 5: Brasil
$ws.Cells.Item(5,2).Value = 1716196   # Casos totales
$ws.Cells.Item(5,3).Value = 41541   # Nuevos casos
$ws.Cells.Item(5,5).Value = 530219   # Recuperados
$ws.Cells.Item(5,7).Value = 1187   # Muertes hoy
$ws.Cells.Item(5,8).Value = 68055   # Muertes

# Row 16: Sudafrica
$ws.Cells.Item(16,5).Value = 114221   # Recuperados
$ws.Cells.Item(16,7).Value = 100   # Muertes hoy
$ws.Cells.Item(16,8).Value = 3602   # Muertes

# Row 22: Colombia
$ws.Cells.Item(22,2).Value = 128638   # Casos totales
$ws.Cells.Item(22,3).Value = 4144   # Nuevos casos
$ws.Cells.Item(22,4).Value = 53634   # Casos activos
$ws.Cells.Item(22,5).Value = 70477   # Recuperados
$ws.Cells.Item(22,7).Value = 168   # Muertes hoy
$ws.Cells.Item(22,8).Value = 4527   # Muertes

# Row 23: Canada
$ws.Cells.Item(23,2).Value = 106415   # Casos totales
$ws.Cells.Item(23,3).Value = 248   # Nuevos casos
$ws.Cells.Item(23,4).Value = 70232   # Casos activos
$ws.Cells.Item(23,5).Value = 27449   # Recuperados
$ws.Cells.Item(23,7).Value = 23   # Muertes hoy
$ws.Cells.Item(23,8).Value = 8734   # Muertes

# Row 25: now "Argentina" (was "China" -- sort order shifted)
$ws.Cells.Item(25,1).Value = "Argentina"
$ws.Cells.Item(25,2).Value = 87030   # Casos totales
$ws.Cells.Item(25,3).Value = 3604   # Nuevos casos
$ws.Cells.Item(25,4).Value = 36502   # Casos activos
$ws.Cells.Item(25,5).Value = 48834   # Recuperados
$ws.Cells.Item(25,7).Value = 50   # Muertes hoy
$ws.Cells.Item(25,8).Value = 1694   # Muertes

# Row 26: now "China" (was "Argentina" -- sort order shifted)
$ws.Cells.Item(26,1).Value = "China"
$ws.Cells.Item(26,2).Value = 83572   # Casos totales
$ws.Cells.Item(26,3).Value = 7   # Nuevos casos
$ws.Cells.Item(26,4).Value = 78548   # Casos activos
$ws.Cells.Item(26,5).Value = 390   # Recuperados
$ws.Cells.Item(26,7).Value = 0   # Muertes hoy
$ws.Cells.Item(26,8).Value = 4634   # Muertes

# Row 44: Panama
$ws.Cells.Item(44,2).Value = 41251   # Casos totales
$ws.Cells.Item(44,3).Value = 960   # Nuevos casos
$ws.Cells.Item(44,4).Value = 19469   # Casos activos
$ws.Cells.Item(44,5).Value = 20963   # Recuperados
$ws.Cells.Item(44,7).Value = 20   # Muertes hoy
$ws.Cells.Item(44,8).Value = 819   # Muertes

# Row 51: now "Nigeria" (was "Rumania" -- sort order shifted)
$ws.Cells.Item(51,1).Value = "Nigeria"
$ws.Cells.Item(51,2).Value = 30249   # Casos totales
$ws.Cells.Item(51,3).Value = 460   # Nuevos casos
$ws.Cells.Item(51,4).Value = 12373   # Casos activos
$ws.Cells.Item(51,5).Value = 17192   # Recuperados
$ws.Cells.Item(51,7).Value = 15   # Muertes hoy
$ws.Cells.Item(51,8).Value = 684   # Muertes

# Row 52: now "Rumania" (was "Armenia" -- sort order shifted)
$ws.Cells.Item(52,1).Value = "Rumania"
$ws.Cells.Item(52,2).Value = 30175   # Casos totales
$ws.Cells.Item(52,3).Value = 555   # Nuevos casos
$ws.Cells.Item(52,4).Value = 20799   # Casos activos
$ws.Cells.Item(52,5).Value = 7559   # Recuperados
$ws.Cells.Item(52,8).Value = 1817   # Muertes

# Row 53: now "Armenia" (was "Nigeria" -- sort order shifted)
$ws.Cells.Item(53,1).Value = "Armenia"
$ws.Cells.Item(53,2).Value = 29820   # Casos totales
$ws.Cells.Item(53,3).Value = 535   # Nuevos casos
$ws.Cells.Item(53,4).Value = 17427   # Casos activos
$ws.Cells.Item(53,5).Value = 11872   # Recuperados
$ws.Cells.Item(53,7).Value = 18   # Muertes hoy
$ws.Cells.Item(53,8).Value = 521   # Muertes

# Row 56: Guatemala
$ws.Cells.Item(56,2).Value = 25411   # Casos totales
$ws.Cells.Item(56,3).Value = 624   # Nuevos casos
$ws.Cells.Item(56,4).Value = 3718   # Casos activos
$ws.Cells.Item(56,5).Value = 20640   # Recuperados
$ws.Cells.Item(56,7).Value = 49   # Muertes hoy
$ws.Cells.Item(56,8).Value = 1053   # Muertes

# Row 57: Ghana
$ws.Cells.Item(57,2).Value = 22822   # Casos totales
$ws.Cells.Item(57,3).Value = 854   # Nuevos casos
$ws.Cells.Item(57,4).Value = 17564   # Casos activos
$ws.Cells.Item(57,5).Value = 5129   # Recuperados

# Row 59: Japon
$ws.Cells.Item(59,2).Value = 20174   # Casos totales
$ws.Cells.Item(59,3).Value = 193   # Nuevos casos
$ws.Cells.Item(59,4).Value = 17331   # Casos activos
$ws.Cells.Item(59,5).Value = 1863   # Recuperados
$ws.Cells.Item(59,7).Value = 2   # Muertes hoy
$ws.Cells.Item(59,8).Value = 980   # Muertes

# Row 69: Chequia
$ws.Cells.Item(69,2).Value = 12814   # Casos totales
$ws.Cells.Item(69,3).Value = 129   # Nuevos casos
$ws.Cells.Item(69,4).Value = 8010   # Casos activos
$ws.Cells.Item(69,5).Value = 4453   # Recuperados

# Row 79: Venezuela
$ws.Cells.Item(79,2).Value = 8008   # Casos totales
$ws.Cells.Item(79,3).Value = 315   # Nuevos casos
$ws.Cells.Item(79,5).Value = 5834   # Recuperados
$ws.Cells.Item(79,7).Value = 3   # Muertes hoy
$ws.Cells.Item(79,8).Value = 74   # Muertes

# Row 141: Uruguay
$ws.Cells.Item(141,2).Value = 974   # Casos totales
$ws.Cells.Item(141,3).Value = 9   # Nuevos casos
$ws.Cells.Item(141,4).Value = 871   # Casos activos
$ws.Cells.Item(141,5).Value = 74   # Recuperados

# Row 209: now "Islas Malvinas" (was "Groenlandia" -- sort order shifted)
$ws.Cells.Item(209,1).Value = "Islas Malvinas"

# Row 210: now "Groenlandia" (was "Islas Malvinas" -- sort order shifted)
$ws.Cells.Item(210,1).Value = "Groenlandia"
